# "Colocando legenda correta nas boxes de 'Industria'."
# Inserts two new columns (entire-column insert, shifting existing data
# right) to hold the "industria" / "industria.es" legend headers that were
# missing from the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column for "industria" goes in before the existing "exportacoes"
# column (currently column E).
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("E1").Value = "industria"
$ws.Range("E1").ColumnWidth = 12.666666666666666

# New column for "industria.es" goes in before the existing
# "exportacoes.es" column (which, after the insert above, is column K).
$ws.Range("K1").EntireColumn.Insert()
$ws.Range("K1").Value = "industria.es"
$ws.Range("K1").ColumnWidth = 10.666666666666666

# Match the workbook's final selection state.
[void]$ws.Range("K2").Select()
